$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Write the table data (headers + 3 data rows) ---
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Age"
$ws.Range("C1").Value = "Location"

$ws.Range("A2").Value = "Jon"
$ws.Range("B2").Value = 21
$ws.Range("C2").Value = "Dallas"

$ws.Range("A3").Value = "Ryan"
$ws.Range("B3").Value = 27
$ws.Range("C3").Value = "Austin"

$ws.Range("A4").Value = "Brady"
$ws.Range("B4").Value = 27
$ws.Range("C4").Value = "Austin"

# Bold the header row before converting to a table, so Excel carries the
# existing direct formatting into the table's header-row style (dxf).
$ws.Range("A1:C1").Font.Bold = $true

# --- Turn the range into a real Excel Table (ListObject) ---
$tbl = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $ws.Range("A1:C4"),
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$tbl.Name = "Table1"
$tbl.TableStyle = "TableStyleLight9"

# --- Column C width ---
$ws.Columns.Item(3).ColumnWidth = 9.6

# --- Page setup ---
$ws.PageSetup.Orientation = [Microsoft.Office.Interop.Excel.XlPageOrientation]::xlPortrait

# --- Selection ---
$ws.Range("H6").Select()
